# Apply the edits described by the commit "Add files via upload":
#   - axes sheet (sheet3) gains three new columns (A_arrow, B_arrow, C_arrow /
#     An (%), Ab (%), Or (%)) inserted before the existing Title/ADAMELLITE
#     column, which is pushed out to column G.
#   - the axes sheet becomes the active/selected sheet (tabSelected) with a
#     new selection of F2, column widths are set on the new columns, and the
#     plys sheet loses its tabSelected flag as a consequence.

$wb = $excel.ActiveWorkbook
$axes = $wb.Worksheets.Item("axes")

# --- axes sheet: shift the old "Title"/"ADAMELLITE" column from D to G and
#     insert the three pairs of new headers/values in D:F ---
$oldTitleHeader = $axes.Range("D1").Value()
$oldTitleValue  = $axes.Range("D2").Value()

$axes.Range("G1").Value = $oldTitleHeader
$axes.Range("G2").Value = $oldTitleValue

$axes.Range("D1").Value = "A_arrow"
$axes.Range("E1").Value = "B_arrow"
$axes.Range("F1").Value = "C_arrow"

$axes.Range("D2").Value = "An (%)"
$axes.Range("E2").Value = "Ab (%)"
$axes.Range("F2").Value = "Or (%)"

# widen the new columns to match the uploaded workbook (closest value the
# engine's column-width rounding grid can represent to the source 14.109375)
$axes.Columns("D:F").ColumnWidth = 13.33

# axes becomes the active sheet/tab, with F2 selected
$axes.Activate() | Out-Null
$axes.Range("F2").Select() | Out-Null
